$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Append a new date entry to the tracked wallet label dates (row 10, column A)
# Force the cell to be stored as plain text (matching the existing date strings)
# rather than being auto-converted into a numeric Excel date serial.
$cell = $ws.Range("A10")
$cell.NumberFormat = "@"
$cell.Value = "2024-10-05"
$cell.Style = "Normal"
